$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2179.4
$ws.Range("J17").Value = 2179.4
$ws.Range("L17").Value = 6538.200000000001
$ws.Range("N17").Value = -6874.200000000001

$ws.Range("H98").Value = 280.33334
$ws.Range("I98").Value = 280.33334
$ws.Range("K98").Value = 280.33334
$ws.Range("M98").Value = 1217.66666

$ws.Range("H99").Value = 318.2
$ws.Range("I99").Value = 318.2
$ws.Range("K99").Value = 954.5999999999999
$ws.Range("M99").Value = 543.4000000000001

$ws.Range("H122").Value = 280.33334
$ws.Range("I122").Value = 280.33334
$ws.Range("K122").Value = 841.0000200000001
$ws.Range("M122").Value = 1608.99998

$ws.Range("H135").Value = 1812.4286
$ws.Range("I135").Value = 1812.4286
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 16311.8574
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -13776.8574

$ws.Range("H138").Value = 6220.8237
$ws.Range("I138").Value = 2073.75
$ws.Range("J138").Value = 9907.111000000001
$ws.Range("K138").Value = 6221.25
$ws.Range("L138").Value = 29721.333
$ws.Range("M138").Value = -1081.25
$ws.Range("N138").Value = -40001.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1882.8823
$ws.Range("I61").Value = 1528.2858
$ws.Range("J61").Value = 3537.6667
$ws.Range("K61").Value = 1528.2858
$ws.Range("L61").Value = 3537.6667
$ws.Range("M61").Value = -1316.2858
$ws.Range("N61").Value = -3961.6667

$ws.Range("H74").Value = 14282600
$ws.Range("I74").Value = 22214600
$ws.Range("K74").Value = 22214600
$ws.Range("M74").Value = -22213726

$ws.Range("H77").Value = 14282600
$ws.Range("I77").Value = 22214600
$ws.Range("K77").Value = 111073000
$ws.Range("M77").Value = -111068632

$ws.Range("H122").Value = 1527.3334
$ws.Range("I122").Value = 1527.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4582.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -2132.0002

$ws.Range("H132").Value = 2032.3096
$ws.Range("I132").Value = 1030.9286
$ws.Range("K132").Value = 3092.7858
$ws.Range("M132").Value = -562.7857999999997

$ws.Range("H136").Value = 1882.8823
$ws.Range("I136").Value = 1528.2858
$ws.Range("J136").Value = 3537.6667
$ws.Range("K136").Value = 4584.857400000001
$ws.Range("L136").Value = 10613.0001
$ws.Range("M136").Value = -2034.857400000001
$ws.Range("N136").Value = -15713.0001

$ws.Range("H139").Value = 51050
$ws.Range("I139").Value = 34650
$ws.Range("J139").Value = 59250
$ws.Range("K139").Value = 34650
$ws.Range("L139").Value = 59250
$ws.Range("M139").Value = -29510
$ws.Range("N139").Value = -69530

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 74890
$ws.Range("J81").Value = 74890
$ws.Range("L81").Value = 74890
$ws.Range("N81").Value = -77012

$ws.Range("H84").Value = 74890
$ws.Range("J84").Value = 74890
$ws.Range("L84").Value = 224670
$ws.Range("N84").Value = -235278

$ws.Range("H134").Value = 3214.6924
$ws.Range("I134").Value = 3214.6924
$ws.Range("K134").Value = 9644.0772
$ws.Range("M134").Value = -7109.0772

$ws.Range("H135").Value = 44418.832
$ws.Range("J135").Value = 44418.832
$ws.Range("L135").Value = 44418.832
$ws.Range("N135").Value = -54558.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1221.75
$ws.Range("I107").Value = 552.2857
$ws.Range("K107").Value = 552.2857
$ws.Range("M107").Value = 1367.7143

$ws.Range("H122").Value = 2302.2144
$ws.Range("I122").Value = 1753
$ws.Range("J122").Value = 3675.25
$ws.Range("K122").Value = 5259
$ws.Range("L122").Value = 11025.75
$ws.Range("M122").Value = -2809
$ws.Range("N122").Value = -15925.75

$ws.Range("H132").Value = 2091.6155
$ws.Range("I132").Value = 1243.8889
$ws.Range("K132").Value = 3731.6667
$ws.Range("M132").Value = -1201.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 96638480
$ws.Range("I4").Value = 56900960
$ws.Range("J4").Value = 160218510
$ws.Range("K4").Value = 170702880
$ws.Range("L4").Value = 480655530
$ws.Range("M4").Value = -170702768
$ws.Range("N4").Value = -480655754

$ws.Range("H131").Value = 1404.7826
$ws.Range("J131").Value = 1554.2941
$ws.Range("L131").Value = 4662.8823
$ws.Range("N131").Value = -14742.8823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3417.8572
$ws.Range("I102").Value = 1637
$ws.Range("K102").Value = 1637
$ws.Range("M102").Value = -15

$ws.Range("H122").Value = 1342.3636
$ws.Range("I122").Value = 1283.75
$ws.Range("J122").Value = 1498.6666
$ws.Range("K122").Value = 3851.25
$ws.Range("L122").Value = 4495.9998
$ws.Range("M122").Value = -1401.25
$ws.Range("N122").Value = -9395.9998

$ws.Range("H126").Value = 3596.6667
$ws.Range("I126").Value = 3596.6667
$ws.Range("K126").Value = 10790.0001
$ws.Range("M126").Value = -8320.000100000001

$ws.Range("H132").Value = 3106.3635
$ws.Range("I132").Value = 1718.25
$ws.Range("K132").Value = 5154.75
$ws.Range("M132").Value = -2624.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3266.4285
$ws.Range("I7").Value = 3227.5
$ws.Range("K7").Value = 3227.5
$ws.Range("M7").Value = -3115.5

$ws.Range("H40").Value = 2395.6365
$ws.Range("I40").Value = 2395.6365
$ws.Range("K40").Value = 2395.6365
$ws.Range("M40").Value = -2259.6365

$ws.Range("H126").Value = 3266.4285
$ws.Range("I126").Value = 3227.5
$ws.Range("K126").Value = 9682.5
$ws.Range("M126").Value = -7212.5

$ws.Range("H132").Value = 4594.4
$ws.Range("I132").Value = 3988.25
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 11964.75
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -9434.75
$ws.Range("N132").Value = -20055.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 5072
$ws.Range("J4").Value = 4250.3335
$ws.Range("L4").Value = 4250.3335
$ws.Range("N4").Value = -4476.3335

$ws.Range("H107").Value = 345.5
$ws.Range("I107").Value = 364
$ws.Range("J107").Value = 216
$ws.Range("K107").Value = 1092
$ws.Range("L107").Value = 648
$ws.Range("M107").Value = 828
$ws.Range("N107").Value = -4488

$ws.Range("H122").Value = 1952.2307
$ws.Range("I122").Value = 1952.2307
$ws.Range("K122").Value = 5856.6921
$ws.Range("M122").Value = -3406.6921

$ws.Range("H132").Value = 3359.2144
$ws.Range("I132").Value = 2714.9333
$ws.Range("K132").Value = 8144.7999
$ws.Range("M132").Value = -5614.7999

$ws.Range("H141").Value = 41863.332
$ws.Range("J141").Value = 41863.332
$ws.Range("L141").Value = 41863.332
$ws.Range("N141").Value = -52223.332
